# Updates cryptos list (Price and Volume(1h) columns) per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to be treated as text so numeric-looking strings
    # (e.g. "302.80") are not coerced into numbers, then drop the
    # temporary text format so the cell keeps its original (default) style.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '42.554.02'
$ws.Range("E2").Value = '  -2.18%  '

Set-TextValue $ws.Range("D3") '2.294.13'
$ws.Range("E3").Value = '  -0.88%  '

$ws.Range("E4").Value = '  -0.06%  '

Set-TextValue $ws.Range("D5") '302.80'
$ws.Range("E5").Value = '  -2.54%  '

Set-TextValue $ws.Range("D6") '98.44'
$ws.Range("E6").Value = '  -6.21%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  -5.66%  '

Set-TextValue $ws.Range("D10") '34.46'
$ws.Range("E10").Value = '  -6.66%  '

$ws.Range("E11").Value = '  -3.27%  '

$ws.Range("E12").Value = '  -0.02%  '

Set-TextValue $ws.Range("D13") '6.73'
$ws.Range("E13").Value = '  -4.42%  '

Set-TextValue $ws.Range("D14") '2.649.36'
$ws.Range("E14").Value = '  -0.85%  '

Set-TextValue $ws.Range("D15") '15.68'
$ws.Range("E15").Value = '  +3.34%  '

Set-TextValue $ws.Range("D16") '2.294.66'
$ws.Range("E16").Value = '  -0.72%  '

Set-TextValue $ws.Range("D17") '0.800'
$ws.Range("E17").Value = '  -1.70%  '

Set-TextValue $ws.Range("D18") '42.475.62'
$ws.Range("E18").Value = '  -2.16%  '

Set-TextValue $ws.Range("D19") '0.0₃0901'
$ws.Range("E19").Value = '  -3.32%  '

Set-TextValue $ws.Range("D20") '11.48'
$ws.Range("E20").Value = '  -6.37%  '

$ws.Range("E21").Value = '  -2.68%  '

Set-TextValue $ws.Range("D22") '67.82'
$ws.Range("E22").Value = '  -0.77%  '

Set-TextValue $ws.Range("D23") '235.21'
$ws.Range("E23").Value = '  -3.22%  '

$ws.Range("E24").Value = '  -3.62%  '

$ws.Range("E25").Value = '  -3.92%  '

$ws.Range("E26").Value = '  -0.10%  '

Set-TextValue $ws.Range("D27") '24.94'
$ws.Range("E27").Value = '  +0.43%  '

Set-TextValue $ws.Range("D29") '34.71'
$ws.Range("E29").Value = '  -6.52%  '

Set-TextValue $ws.Range("D30") '9.17'
$ws.Range("E30").Value = '  -5.24%  '

Set-TextValue $ws.Range("D31") '163.42'
$ws.Range("E31").Value = '  -1.91%  '

Set-TextValue $ws.Range("D32") '1.00'
$ws.Range("E32").Value = '  +0.00%  '

Set-TextValue $ws.Range("D33") '5.01'
$ws.Range("E33").Value = '  -5.70%  '

Set-TextValue $ws.Range("D34") '4.58'
$ws.Range("E34").Value = '  -0.51%  '

$ws.Range("E35").Value = '  -5.03%  '

$ws.Range("E36").Value = '  -4.87%  '

Set-TextValue $ws.Range("D37") '16.92'
$ws.Range("E37").Value = '  -8.09%  '

Set-TextValue $ws.Range("D38") '2.88'
$ws.Range("E38").Value = '  -6.23%  '

Set-TextValue $ws.Range("D39") '1.80'
$ws.Range("E39").Value = '  -4.82%  '

$ws.Range("E40").Value = '  -5.91%  '

$ws.Range("E41").Value = '  -4.06%  '

Set-TextValue $ws.Range("D42") '2.43'
$ws.Range("E42").Value = '  -11.29%  '

Set-TextValue $ws.Range("D43") '1.979.46'
$ws.Range("E43").Value = '  -0.94%  '

Set-TextValue $ws.Range("D44") '0.0279'
$ws.Range("E44").Value = '  -4.89%  '

Set-TextValue $ws.Range("D45") '18.56'
$ws.Range("E45").Value = '  -2.92%  '

Set-TextValue $ws.Range("D46") '10.17'
$ws.Range("E46").Value = '  +1.29%  '

Set-TextValue $ws.Range("D47") '2.91'
$ws.Range("E47").Value = '  -7.77%  '

Set-TextValue $ws.Range("D48") '55.44'
$ws.Range("E48").Value = '  -2.75%  '

$ws.Range("E49").Value = '  -3.68%  '

Set-TextValue $ws.Range("D50") '2.517.52'
$ws.Range("E50").Value = '  -0.93%  '

Set-TextValue $ws.Range("D51") '4.68'
$ws.Range("E51").Value = '  -0.91%  '
